$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.407.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.100.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5218"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4557"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.20"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +13.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08890"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.178"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.088.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.796"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.014"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001146"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06617"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.294"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.473.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.358"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.336.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.514"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.206"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1066"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.651"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.393"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.934"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.859"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02575"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2313"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6872"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.249"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.319"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6394"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.655"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.246"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "83.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000338"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.99%  "
